$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: yessirski/hahaha123 -> admin/admin
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "admin"

# Update row 3: admin/admin -> 123/1 (keep as text, not numbers)
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "123"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"

# Remove old row 4 (123/123) entirely so dimension becomes A1:B3
$ws.Rows.Item(4).Delete()
